$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.445.13"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "'1.852.29"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'240.82"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'0.6304"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07683"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "'0.2943"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "'24.64"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'1.851.52"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.00001097"
$ws.Range("E13").Value = "  +7.17%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.027"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'0.6814"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'83.61"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'2.105.70"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'6.169"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'29.457.87"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'229.58"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'7.458"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'157.04"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'0.1389"
$ws.Range("D27").Value = "'8.405"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'0.05698"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'4.055"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "'1.851"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'1.163"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'0.7050"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "'2.586"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'2.783"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "'1.220.00"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "'6.523"
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("D42").Value = "'0.9086"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'2.014.44"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'101.81"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'66.51"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'0.00000000119"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'7.134"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").Value = "'0.4021"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.006"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.686"
$ws.Range("E51").Value = "  -1.27%  "
